$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G header: "curls"
$ws.Range("G1").Value = "curls"

# Existing rows 2-37 get a 0 in the new "curls" column
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}

# New row 38: date 2021-02-25 (serial 44252), with values
$ws.Cells.Item(38, 1).Value = 44252
$ws.Cells.Item(38, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(38, 2).Value = 0
$ws.Cells.Item(38, 3).Value = 60
$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 32

# Update the selection to match the saved view state
$ws.Range("L23").Select()
